# "La til klimatall 2021" - re-arranges the municipality rows so that the
# three "newly added 2021" rows (Larvik, Holmestrand, Horten) move below the
# rest of the list, leaving three blank spacer rows behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the six rows that stay in the list up by three rows (16-21) ---
$ws.Range("A16").Value = "Siljan"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "Siljan"

$ws.Range("A17").Value = "Færder"
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "Færder"

$ws.Range("A18").Value = "Kragerø"
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = "Kragerø"

$ws.Range("A19").Value = "Bamble"
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = "Bamble"

$ws.Range("A20").Value = "Notodden"
$ws.Range("B20").Value = 16
$ws.Range("C20").Value = "Notodden"

$ws.Range("A21").Value = "Skien"
$ws.Range("B21").Value = 19
$ws.Range("C21").Value = "Skien"

# --- Three blank spacer rows (22-24) ---
$ws.Range("A22:C24").ClearContents()

# --- The three moved rows land at the bottom of the list (25-27) ---
$ws.Range("A25").Value = "Larvik"
$ws.Range("B25").Value = 15
$ws.Range("C25").Value = "Larvik"

$ws.Range("A26").Value = "Holmestrand"
$ws.Range("B26").Value = 12
$ws.Range("C26").Value = "Holmestrand"

$ws.Range("A27").Value = "Horten"
$ws.Range("B27").Value = 20
$ws.Range("C27").Value = "Horten"

# --- Column B widened slightly ---
$ws.Columns("B").ColumnWidth = 10

# --- Active cell / selection moved to D29 ---
$ws.Range("D29").Select()
